$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new rows 30-32 with new programs (this introduces new shared strings
# in the order: Git 2.9+, Em progresso, Git Hub, Atom)
$ws.Range("A30").Value = "Git 2.9+"
$ws.Range("B30").Value = "Em progresso"
$ws.Range("A31").Value = "Git Hub"
$ws.Range("B31").Value = "Pendente"
$ws.Range("A32").Value = "Atom"
$ws.Range("B32").Value = "Pendente"

# Update existing rows' status (reusing "Em progresso", adding new "Resolvido")
$ws.Range("B26").Value = "Em progresso"
$ws.Range("B29").Value = "Resolvido"

# Update the visible scroll/selection state to match the new view
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 21
$ws.Range("B29").Select()
